$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes a string so it survives round-trip as text even when the
# content looks numeric (e.g. "0.780", "19.20") -- Excel normally
# auto-converts such strings to numbers and drops trailing zeros.
function Set-TextValue($Cell, $Text) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2,4) "37.396.38"
$ws.Cells.Item(2,5).Value = "  +3.34%  "
Set-TextValue $ws.Cells.Item(3,4) "2.061.66"
$ws.Cells.Item(4,5).Value = "  +0.07%  "
Set-TextValue $ws.Cells.Item(5,4) "237.26"
$ws.Cells.Item(5,5).Value = "  +3.14%  "
Set-TextValue $ws.Cells.Item(6,4) "0.617"
$ws.Cells.Item(6,5).Value = "  +3.55%  "
Set-TextValue $ws.Cells.Item(7,4) "58.19"
$ws.Cells.Item(7,5).Value = "  +8.49%  "
$ws.Cells.Item(8,5).Value = "  -0.02%  "
$ws.Cells.Item(9,5).Value = "  +4.87%  "
Set-TextValue $ws.Cells.Item(10,4) "57.58"
$ws.Cells.Item(10,5).Value = "  +0.91%  "
Set-TextValue $ws.Cells.Item(11,4) "0.0762"
$ws.Cells.Item(11,5).Value = "  +2.46%  "
$ws.Cells.Item(12,5).Value = "  +3.86%  "
Set-TextValue $ws.Cells.Item(13,4) "2.364.13"
$ws.Cells.Item(13,5).Value = "  +4.69%  "
Set-TextValue $ws.Cells.Item(14,4) "14.43"
$ws.Cells.Item(14,5).Value = "  +4.97%  "
Set-TextValue $ws.Cells.Item(15,4) "21.18"
$ws.Cells.Item(15,5).Value = "  +7.33%  "
Set-TextValue $ws.Cells.Item(16,4) "0.780"
$ws.Cells.Item(16,5).Value = "  +4.65%  "
Set-TextValue $ws.Cells.Item(17,4) "5.19"
$ws.Cells.Item(17,5).Value = "  +4.07%  "
Set-TextValue $ws.Cells.Item(18,4) "2.065.27"
$ws.Cells.Item(18,5).Value = "  +4.94%  "
Set-TextValue $ws.Cells.Item(19,4) "37.586.24"
$ws.Cells.Item(19,5).Value = "  +3.88%  "
Set-TextValue $ws.Cells.Item(20,4) "6.14"
$ws.Cells.Item(20,5).Value = "  +21.73%  "
Set-TextValue $ws.Cells.Item(21,4) "69.05"
$ws.Cells.Item(21,5).Value = "  +2.76%  "
Set-TextValue $ws.Cells.Item(22,4) "0.0₃0814"
$ws.Cells.Item(22,5).Value = "  +1.76%  "
Set-TextValue $ws.Cells.Item(23,4) "225.55"
$ws.Cells.Item(23,5).Value = "  +2.54%  "
$ws.Cells.Item(24,5).Value = "  +0.00%  "
Set-TextValue $ws.Cells.Item(25,4) "2.46"
$ws.Cells.Item(25,5).Value = "  +6.02%  "
$ws.Cells.Item(26,5).Value = "  +2.68%  "
Set-TextValue $ws.Cells.Item(27,4) "163.85"
$ws.Cells.Item(27,5).Value = "  +2.63%  "
$ws.Cells.Item(28,2).Value = "ImmutableX"
$ws.Cells.Item(28,3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Cells.Item(28,4) "1.46"
$ws.Cells.Item(28,5).Value = "  +11.57%  "
$ws.Cells.Item(29,2).Value = "Cosmos"
$ws.Cells.Item(29,3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Cells.Item(29,4) "8.88"
$ws.Cells.Item(29,5).Value = "  +4.95%  "
Set-TextValue $ws.Cells.Item(30,4) "19.20"
$ws.Cells.Item(30,5).Value = "  +3.38%  "
$ws.Cells.Item(31,5).Value = "  +4.47%  "
$ws.Cells.Item(32,5).Value = "  +3.47%  "
$ws.Cells.Item(33,5).Value = "  +4.51%  "
$ws.Cells.Item(34,5).Value = "  +4.89%  "
$ws.Cells.Item(35,5).Value = "  +13.12%  "
$ws.Cells.Item(36,5).Value = "  +6.17%  "
$ws.Cells.Item(37,5).Value = "  +0.05%  "
Set-TextValue $ws.Cells.Item(38,4) "3.37"
$ws.Cells.Item(38,5).Value = "  +5.33%  "
$ws.Cells.Item(39,5).Value = "  +0.62%  "
Set-TextValue $ws.Cells.Item(40,4) "5.88"
$ws.Cells.Item(40,5).Value = "  +14.21%  "
$ws.Cells.Item(41,5).Value = "  +11.89%  "
$ws.Cells.Item(42,2).Value = "HuobiToken"
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Cells.Item(42,4) "2.97"
$ws.Cells.Item(42,5).Value = "  -2.19%  "
$ws.Cells.Item(43,2).Value = "FTXToken"
$ws.Cells.Item(43,3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Cells.Item(43,4) "4.47"
$ws.Cells.Item(43,5).Value = "  +23.24%  "
Set-TextValue $ws.Cells.Item(44,4) "97.50"
$ws.Cells.Item(44,5).Value = "  +12.01%  "
Set-TextValue $ws.Cells.Item(45,4) "1.481.50"
$ws.Cells.Item(45,5).Value = "  +3.95%  "
Set-TextValue $ws.Cells.Item(46,4) "0.0210"
$ws.Cells.Item(46,5).Value = "  +6.03%  "
$ws.Cells.Item(47,5).Value = "  +7.60%  "
$ws.Cells.Item(48,5).Value = "  +9.12%  "
Set-TextValue $ws.Cells.Item(49,4) "1.03"
$ws.Cells.Item(49,5).Value = "  +4.58%  "
Set-TextValue $ws.Cells.Item(50,4) "7.21"
$ws.Cells.Item(50,5).Value = "  +7.46%  "
$ws.Cells.Item(51,5).Value = "  +2.94%  "
